$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row: rename / re-key the columns and add the new "key" column in A.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "key"
$ws.Range("B1").Value = "junction name"
$ws.Range("C1").Value = "numOfLanes"

# ---------------------------------------------------------------------------
# Junction rows: new names + corrected / re-ordered compass directions.
# Column layout stays A=key, B=name, C=numOfLanes, D..G = direction labels.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "A"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "east"
$ws.Range("E2").Value = "west"
$ws.Range("F2").Value = "north"
$ws.Range("G2").Value = "south"

$ws.Range("B3").Value = "Morasha"
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "west"
$ws.Range("E3").Value = "east"
$ws.Range("F3").Value = "north"
$ws.Range("G3").Value = "south"

$ws.Range("B4").Value = "Yarkon"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "east"
$ws.Range("E4").Value = "west"
$ws.Range("F4").Value = "north"
$ws.Range("G4").Value = "south"

$ws.Range("B5").Value = "C"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "south"
$ws.Range("E5").Value = "north"
$ws.Range("F5").Value = "ease"
$ws.Range("G5").Value = "west"

$ws.Range("B6").Value = "BarIlan"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = "north"
$ws.Range("E6").Value = "west"
$ws.Range("F6").Value = "south"
$ws.Range("G6").Value = "east"

$ws.Range("B7").Value = "AlufSafe"
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = "east"
$ws.Range("E7").Value = "west"
$ws.Range("F7").Value = "north"
$ws.Range("G7").Value = "south"

# ---------------------------------------------------------------------------
# Column widths (A1 notation header columns through I), mirroring how a
# user would drag/resize columns so the new data is readable.
# Excel's ColumnWidth is expressed in characters and gets rounded to a
# whole-pixel grid on write, so we back-solve the character width that
# reproduces the desired on-disk width as closely as possible.
# ---------------------------------------------------------------------------
function Set-ColWidth($colIndex, [double]$targetChars) {
    $mdw = 7.0
    $px = [Math]::Round($targetChars * $mdw)
    $chars = ($px - 5) / $mdw
    $ws.Columns.Item($colIndex).ColumnWidth = $chars
}

Set-ColWidth 1 13.5
Set-ColWidth 2 18.796875
Set-ColWidth 3 20.296875
Set-ColWidth 4 16.5
Set-ColWidth 5 17.3984375
Set-ColWidth 6 12.3984375
Set-ColWidth 7 18.5
Set-ColWidth 8 17.69921875
Set-ColWidth 9 16.59765625

# ---------------------------------------------------------------------------
# Selection follows the last-edited cell.
# ---------------------------------------------------------------------------
$ws.Range("B6").Select()

Write-Host "Junction_Info updated"
